$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6460902094841003
$ws.Range("B1").Value = 0.6371238827705383
$ws.Range("C1").Value = 0.6677056550979614
$ws.Range("D1").Value = 0.8804638385772705
$ws.Range("E1").Value = 0.8508303165435791
